# Update the "Förändrad" (column C) date value for rows 2-27 from 45327 to 45328
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C27").Value = 45328

# Row 27 loses its explicit custom row height (reverts to default sheet height)
$ws.Rows.Item(27).AutoFit()

# Row 28 (A 4345-2024) is removed entirely
$ws.Rows.Item(28).Delete()
